# Revert "User data 3.0"
#
# On the "Data-wide-value" sheet the previous edit had inserted a new
# column B ("budget-type", all rows literally "budget") in front of the
# five wide year columns (2012..2016), pushing the real data from
# B:F out to C:G. Reverting means removing that inserted column again:
# delete column B and shift everything back to the left, so the sheet
# goes back to id | 2012 | 2013 | 2014 | 2015 | 2016 (A:F).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data-wide-value")

$ws.Range("B1:B112").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)
